$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 242 (pushes the existing 242:308 block down to 244:310,
# and the sheet dimension grows from R308 to R310). Excel's native row-insert
# carries the date number-format down from the row above for column D, which
# matches the workbook's existing styling for that column.
$ws.Rows("242:243").Insert()

# Row 242 - new weekly entry, "Primera" quality
$ws.Cells.Item(242, 1).Value  = 5
$ws.Cells.Item(242, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(242, 3).Value  = "Maule"
$ws.Cells.Item(242, 4).Value  = 44754
$ws.Cells.Item(242, 5).Value  = 7
$ws.Cells.Item(242, 6).Value  = 100112006
$ws.Cells.Item(242, 7).Value  = "Repollo"
$ws.Cells.Item(242, 8).Value  = "Crespo record"
$ws.Cells.Item(242, 9).Value  = "Primera"
$ws.Cells.Item(242, 10).Value = 2000
$ws.Cells.Item(242, 11).Value = 1200
$ws.Cells.Item(242, 12).Value = 1200
$ws.Cells.Item(242, 13).Value = 1200
$ws.Cells.Item(242, 14).Value = "`$/unidad"
$ws.Cells.Item(242, 15).Value = "Región del Maule"
$ws.Cells.Item(242, 16).Value = 1200
$ws.Cells.Item(242, 17).Value = 1
$ws.Cells.Item(242, 18).Value = "Hortaliza"

# Row 243 - new weekly entry, "Segunda" quality
$ws.Cells.Item(243, 1).Value  = 5
$ws.Cells.Item(243, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(243, 3).Value  = "Maule"
$ws.Cells.Item(243, 4).Value  = 44754
$ws.Cells.Item(243, 5).Value  = 7
$ws.Cells.Item(243, 6).Value  = 100112006
$ws.Cells.Item(243, 7).Value  = "Repollo"
$ws.Cells.Item(243, 8).Value  = "Crespo record"
$ws.Cells.Item(243, 9).Value  = "Segunda"
$ws.Cells.Item(243, 10).Value = 2000
$ws.Cells.Item(243, 11).Value = 1000
$ws.Cells.Item(243, 12).Value = 1000
$ws.Cells.Item(243, 13).Value = 1000
$ws.Cells.Item(243, 14).Value = "`$/unidad"
$ws.Cells.Item(243, 15).Value = "Región del Maule"
$ws.Cells.Item(243, 16).Value = 1000
$ws.Cells.Item(243, 17).Value = 1
$ws.Cells.Item(243, 18).Value = "Hortaliza"
